$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 90

$ws.Range("B${row}:D${row}").NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = 1581552000
$ws.Cells.Item($row, 2).Value = "2020-02-13"
$ws.Cells.Item($row, 3).Value = "0213"
$ws.Cells.Item($row, 4).Value = "MTAG"
$ws.Cells.Item($row, 5).Value = 0.5
$ws.Cells.Item($row, 6).Value = 0.505
$ws.Cells.Item($row, 7).Value = 0.495
$ws.Cells.Item($row, 8).Value = 0.5
$ws.Cells.Item($row, 9).Value = 4920500

$ws.Range("B${row}:D${row}").ClearFormats()
